# The document contains several inline OMath (equation) fields used as
# arrow/plus-minus glyphs between treatment steps, e.g. "Chemo [arrow] Surgery".
# Most of them already store a single run with the Unicode glyph directly
# (styled with <m:sty m:val="p"/> so Word renders it upright, not italic).
#
# Two of them - in the "Superficial (T1)" and "Localized (T1b/T2)" bullet
# points under "4 Treatment Plan" - instead spell out the word "Rightarrow"
# across ten separate <m:r><m:t>…</m:t></m:r> runs (R-i-g-h-t-a-r-r-o-w).
# This script collapses each of those two equations down to a single,
# properly-styled run containing the actual "RIGHTWARDS DOUBLE ARROW"
# glyph (U+21D2), matching the style/glyph already used by the other,
# similar "implies next treatment step" equations in the document.

$d = $word.ActiveDocument

# U+21D2 = RIGHTWARDS DOUBLE ARROW ("⇒")
$arrow = [char]0x21D2
$replacementXml = "<m:oMath><m:r><m:rPr><m:sty m:val=""p""/></m:rPr><m:t>$arrow</m:t></m:r></m:oMath>"

# Plain ASCII "R" stored inside an <m:t> run renders (via Range.Text) as
# the math-italic capital R (U+1D445), since Word italicizes bare Latin
# letters inside math zones by default. The ten-run spelled-out
# "Rightarrow" equations therefore show up as a 19-"character" Range.Text
# (10 real glyphs + 9 empty run-boundary slots) starting with that italic
# "R", while every other equation in the document is a single symbol
# (Range.Text length 1). That combination safely identifies only the two
# equations that need to change.
$italicCapR = [char]0x1D445

$targets = @()
for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    $t = $om.Range.Text
    if ($t.Length -eq 19 -and $t.Substring(0, 1) -eq $italicCapR) {
        $targets += $i
    }
}

# Replace from the last match backwards so earlier indices/ranges in the
# OMaths collection stay valid while we edit later ones.
for ($k = $targets.Count - 1; $k -ge 0; $k--) {
    $idx = $targets[$k]
    $om = $d.OMaths.Item($idx)
    $om.Range.InsertXML($replacementXml)
}

Write-Host "Replaced $($targets.Count) spelled-out 'Rightarrow' equations with the arrow glyph."
